# Scheduled-runner refresh of market-price-derived columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-job
# Moogle Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values below
# mirror the new snapshot pulled from the price API; columns untouched by the
# refresh (A-G and any cell not listed) are left exactly as-is.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 619.75
$ws.Range("I2").Value = 147.11111
$ws.Range("K2").Value = 147.11111
$ws.Range("M2").Value = -34.11111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1027.5
$ws.Range("J5").Value = 97.666664
$ws.Range("L5").Value = 97.666664
$ws.Range("N5").Value = -327.666664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1606.2433
$ws.Range("J17").Value = 1606.2433
$ws.Range("L17").Value = 4818.7299
$ws.Range("N17").Value = -5154.7299

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 237.4
$ws.Range("J39").Value = 482
$ws.Range("L39").Value = 1446
$ws.Range("N39").Value = -2038

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 1116.1818
$ws.Range("I42").Value = 189.5
$ws.Range("J42").Value = 1322.1111
$ws.Range("K42").Value = 568.5
$ws.Range("L42").Value = 3966.3333
$ws.Range("M42").Value = -338.5
$ws.Range("N42").Value = -4426.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 19499.75
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2823.8
$ws.Range("I70").Value = 2330
$ws.Range("J70").Value = 3035.4285
$ws.Range("K70").Value = 6990
$ws.Range("L70").Value = 9106.2855
$ws.Range("M70").Value = -6720
$ws.Range("N70").Value = -9646.2855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 19499.75
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2823.8
$ws.Range("I73").Value = 2330
$ws.Range("J73").Value = 3035.4285
$ws.Range("K73").Value = 6990
$ws.Range("L73").Value = 9106.2855
$ws.Range("M73").Value = -6054
$ws.Range("N73").Value = -10978.2855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7217.6665
$ws.Range("I116").Value = 6731.3335
$ws.Range("K116").Value = 6731.3335
$ws.Range("M116").Value = -3289.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2621.4167
$ws.Range("I132").Value = 2621.4167
$ws.Range("K132").Value = 7864.250100000001
$ws.Range("M132").Value = -5334.250100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 5054.8125
$ws.Range("I135").Value = 1722.3334
$ws.Range("K135").Value = 15501.0006
$ws.Range("M135").Value = -12966.0006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2703.9333
$ws.Range("I137").Value = 1767
$ws.Range("K137").Value = 5301
$ws.Range("M137").Value = -2751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1838.2307
$ws.Range("I45").Value = 1144.85
$ws.Range("J45").Value = 4149.5
$ws.Range("K45").Value = 1144.85
$ws.Range("L45").Value = 4149.5
$ws.Range("M45").Value = -767.8499999999999
$ws.Range("N45").Value = -4903.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3700.5186
$ws.Range("I122").Value = 3554.6667
$ws.Range("K122").Value = 10664.0001
$ws.Range("M122").Value = -8214.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 52800.855
$ws.Range("J21").Value = 52800.855
$ws.Range("L21").Value = 52800.855
$ws.Range("N21").Value = -53272.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 27163.715
$ws.Range("J54").Value = 43747
$ws.Range("L54").Value = 43747
$ws.Range("N54").Value = -44715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1166.0952
$ws.Range("I86").Value = 1044.4667
$ws.Range("K86").Value = 1044.4667
$ws.Range("M86").Value = 78.53330000000005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1166.0952
$ws.Range("I89").Value = 1044.4667
$ws.Range("K89").Value = 5222.3335
$ws.Range("M89").Value = 393.6665000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3766.9
$ws.Range("I99").Value = 2627.25
$ws.Range("J99").Value = 5476.375
$ws.Range("K99").Value = 2627.25
$ws.Range("L99").Value = 5476.375
$ws.Range("M99").Value = -1129.25
$ws.Range("N99").Value = -8472.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 156392.86
$ws.Range("J106").Value = 156392.86
$ws.Range("L106").Value = 156392.86
$ws.Range("N106").Value = -158916.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H118").Value = 199950
$ws.Range("J118").Value = 199950
$ws.Range("L118").Value = 199950
$ws.Range("N118").Value = -203264

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 617.3143
$ws.Range("J7").Value = 483.4
$ws.Range("L7").Value = 483.4
$ws.Range("N7").Value = -709.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6429.963
$ws.Range("I58").Value = 2674.3157
$ws.Range("J58").Value = 15349.625
$ws.Range("K58").Value = 2674.3157
$ws.Range("L58").Value = 15349.625
$ws.Range("M58").Value = -2471.3157
$ws.Range("N58").Value = -15755.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2329.6
$ws.Range("I99").Value = 1943.3125
$ws.Range("K99").Value = 1943.3125
$ws.Range("M99").Value = -445.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2329.6
$ws.Range("I126").Value = 1943.3125
$ws.Range("K126").Value = 5829.9375
$ws.Range("M126").Value = -3359.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4270.8857
$ws.Range("I132").Value = 2981.074
$ws.Range("K132").Value = 8943.222
$ws.Range("M132").Value = -6413.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6429.963
$ws.Range("I136").Value = 2674.3157
$ws.Range("J136").Value = 15349.625
$ws.Range("K136").Value = 8022.9471
$ws.Range("L136").Value = 46048.875
$ws.Range("M136").Value = -5472.9471
$ws.Range("N136").Value = -51148.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 449
$ws.Range("I5").Value = 449
$ws.Range("K5").Value = 1347
$ws.Range("M5").Value = -1235

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 16169
$ws.Range("J112").Value = 18585.8
$ws.Range("L112").Value = 55757.39999999999
$ws.Range("N112").Value = -57973.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1680.875
$ws.Range("I113").Value = 1379
$ws.Range("J113").Value = 1724
$ws.Range("K113").Value = 4137
$ws.Range("L113").Value = 5172
$ws.Range("M113").Value = -1967
$ws.Range("N113").Value = -9512

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 449
$ws.Range("I135").Value = 449
$ws.Range("K135").Value = 4041
$ws.Range("M135").Value = -1506

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6745.9355
$ws.Range("I80").Value = 6366.75
$ws.Range("K80").Value = 6366.75
$ws.Range("M80").Value = -5368.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6745.9355
$ws.Range("I83").Value = 6366.75
$ws.Range("K83").Value = 31833.75
$ws.Range("M83").Value = -26841.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 109876.86
$ws.Range("J100").Value = 109876.86
$ws.Range("L100").Value = 109876.86
$ws.Range("N100").Value = -112040.86

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3076.353
$ws.Range("I102").Value = 1881.7142
$ws.Range("K102").Value = 1881.7142
$ws.Range("M102").Value = -259.7141999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11759.8
$ws.Range("I122").Value = 10856.857
$ws.Range("K122").Value = 32570.571
$ws.Range("M122").Value = -30120.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 578.1905
$ws.Range("I16").Value = 507.35
$ws.Range("K16").Value = 507.35
$ws.Range("M16").Value = -337.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1472.4375
$ws.Range("I22").Value = 656
$ws.Range("K22").Value = 656
$ws.Range("M22").Value = -361

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1472.4375
$ws.Range("I27").Value = 656
$ws.Range("K27").Value = 656
$ws.Range("M27").Value = -549

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3400
$ws.Range("I46").Value = 1520
$ws.Range("J46").Value = 4575
$ws.Range("K46").Value = 1520
$ws.Range("L46").Value = 4575
$ws.Range("M46").Value = -1332
$ws.Range("N46").Value = -4951

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1635.3182
$ws.Range("I93").Value = 1285.2142
$ws.Range("J93").Value = 2248
$ws.Range("K93").Value = 1285.2142
$ws.Range("L93").Value = 2248
$ws.Range("M93").Value = -37.21419999999989
$ws.Range("N93").Value = -4744

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 159964.14
$ws.Range("J102").Value = 178291.5
$ws.Range("L102").Value = 178291.5
$ws.Range("N102").Value = -184781.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3108.7222
$ws.Range("I132").Value = 1197.5333
$ws.Range("K132").Value = 3592.5999
$ws.Range("M132").Value = -1062.5999
